$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in new "hw11" (column L) grades for rows 2-10
$ws.Range("L2").Value = 97
$ws.Range("L3").Value = 98
$ws.Range("L4").Value = 99
$ws.Range("L5").Value = 97
$ws.Range("L6").Value = 89
$ws.Range("L7").Value = 100
$ws.Range("L8").Value = 100
$ws.Range("L9").Value = 100
$ws.Range("L10").Value = 100

# Update the active selection to reflect where the user left off editing
$ws.Range("L11").Select()
